# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45186 (2023-09-17) to 45188 (2023-09-19).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Determine the last used row reliably (column A holds the record id for
# every data row), data starts on row 2 (row 1 is the header row).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -ne $null) {
        $cell.Value = 45188
    }
}
